# Auto-generated Excel COM-interop script updating the cryptos worksheet
# per the commit "Updated cryptos list on Tue Jul 16 19:25:36 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.930.55"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.90%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.462.40"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.39%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.35%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.00"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.35%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.13%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.460.56"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.05%  "

# Row 9
$ws.Range("E9").Value = "  +8.50%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.33"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.73%  "

# Row 11
$ws.Range("E11").Value = "  +3.52%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.441"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.96%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.051.10"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.75%  "

# Row 14
$ws.Range("E14").Value = "  -2.94%  "

# Row 15
$ws.Range("E15").Value = "  +4.53%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.76"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +5.36%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.897.40"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.84%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.462.64"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.72%  "

# Row 19
$ws.Range("E19").Value = "  -0.43%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.35"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.84%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "389.25"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.08%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.21"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.19%  "

# Row 23
$ws.Range("E23").Value = "  +1.87%  "

# Row 24
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.01"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.46%  "

# Row 25
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.04"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.73%  "

# Row 26
$ws.Range("E26").Value = "  +16.94%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.52"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.15%  "

# Row 28
$ws.Range("E28").Value = "  +0.31%  "

# Row 29
$ws.Range("E29").Value = "  +0.00%  "

# Row 30
$ws.Range("E30").Value = "  +8.55%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.45"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +7.43%  "

# Row 32
$ws.Range("E32").Value = "  -0.03%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.67"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.75%  "

# Row 34
$ws.Range("E34").Value = "  -0.48%  "

# Row 35
$ws.Range("E35").Value = "  +0.17%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.10"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.59%  "

# Row 37
$ws.Range("E37").Value = "  +0.51%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.94"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.28%  "

# Row 39
$ws.Range("E39").Value = "  +1.27%  "

# Row 40
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.959.02"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.32%  "

# Row 41
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0768"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.06%  "

# Row 42
$ws.Range("E42").Value = "  -1.75%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.58"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.48%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.87"
$ws.Range("D44").ClearFormats()

# Row 45
$ws.Range("E45").Value = "  -2.42%  "

# Row 46
$ws.Range("E46").Value = "  +1.23%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.11"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +7.05%  "

# Row 48
$ws.Range("E48").Value = "  +1.96%  "

# Row 49
$ws.Range("E49").Value = "  +10.96%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.875"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.55%  "

# Row 51
$ws.Range("E51").Value = "  +3.69%  "
